# "Generate Report for Handoff" — refresh the Latest Handoff / Latest HO
# Xliff Generate Date timestamps for the d101eedf-e7e5-42e6-9713-4ccd5eb4805e
# row (row 5 of each table) across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-10-21 00:08:13"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-10-21 00:08:02"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-10-21 00:08:13"
